# Fix typo in peak map: "18:2 trans 9,12" -> "18:2 trans 9, 12"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A10").Value = "18:2 trans 9, 12"
